$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before "batsman" (currently column D) for ownTeam / oppTeam
$ws.Columns("D:E").Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Row 2 (existing Sharjah match vs Mumbai) - fill new team columns
$ws.Range("D2").Value = "Chennai Super Kings"
$ws.Range("E2").Value = "Mumbai Indians"

# New row 3: Dubai match vs Royal Challengers Bangalore
# Force the numeric-looking text cells (G3:K3) to stay stored as text, same as
# the rest of the sheet (numberStoredAsText ignored error covers the whole range).
$ws.Range("G3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 10 2020"
$ws.Range("C3").Value = "RCB won by 37 runs"
$ws.Range("D3").Value = "Chennai Super Kings"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Shardul Thakur "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "100.00"
